$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the two customer-number cells.
$ws.Range("P2").Value = 509514334
$ws.Range("P3").Value = 504798215

# Move the active selection to P4 (matches the saved cursor position).
$ws.Range("P4").Select()
